$d = $word.ActiveDocument

# --- Topic block 1 ---
$d.Paragraphs(1).Range.Text = 'Topic: GitHub Contribution Graph Manipulation'
$d.Paragraphs(3).Range.Text = 'The GitHub contribution graph visually represents a user''s activity on the platform, showing the number of commits made on specific days. This tutorial focuses on manipulating this graph to display a desired pattern or "art." Understanding how GitHub tracks commits is essential; each commit is associated with a timestamp. By controlling these timestamps, you can effectively alter the appearance of your contribution graph. This technique involves automating commits and linking them to past dates, making them appear as if they were made historically. The motivation behind the tutorial is primarily educational, to highlight how easily one can create an appearance of activity on the internet that may not reflect reality.'
$d.Paragraphs(5).Range.Text = '*   The contribution graph visually represents a user''s commit history.'
$d.Paragraphs(6).Range.Text = '*   Manipulating this graph involves controlling commit timestamps.'
$d.Paragraphs(7).Range.Text = '*   The primary goal is to demonstrate the ability to fake contributions.'
$d.Paragraphs(8).Range.Text = '*   Automation is key for creating a large number of commits over time.'
$d.Paragraphs(9).Range.Text = '*   This aims to show how easy it is to create an illusion of activity.'
$d.Paragraphs(11).Range.Delete()

# --- Topic block 2 ---
$d.Paragraphs(11).Range.Text = 'Topic: Timestamps in GitHub Commits'
$d.Paragraphs(13).Range.Text = 'GitHub stores each commit with a timestamp, which is a crucial element for understanding the process of manipulating the contribution graph. This timestamp records the date and time when the commit was made. The contribution graph uses these timestamps to determine where to place the contributions visually. By manipulating the timestamps associated with commits, you can influence how your contributions are displayed on your profile. This tutorial leverages timestamps to make it appear that you made commits in the past or on specific dates. The core idea is to mimic how GitHub itself stores this data, thus fooling the system and making you appear active on certain days of the past, regardless of the actual commit date.'
$d.Paragraphs(15).Range.Text = '*   Each commit is paired with a timestamp by GitHub.'
$d.Paragraphs(16).Range.Text = '*   The timestamp dictates the position on the contribution graph.'
$d.Paragraphs(17).Range.Text = '*   Manipulating timestamps allows for control over the graph''s appearance.'
$d.Paragraphs(18).Range.Text = '*   The technique involves associating commits with past dates.'
$d.Paragraphs(19).Range.Text = '*   The focus here is on mirroring GitHub''s method of storing data.'
$d.Paragraphs(21).Range.Delete()

# --- Topic block 3 ---
$d.Paragraphs(21).Range.Text = 'Topic: Automating Commits'
$d.Paragraphs(23).Range.Text = 'Automating commits is central to creating the desired patterns on the GitHub contribution graph. This tutorial outlines how to generate numerous commits and associate them with specific dates, creating an illusion of historical activity. Automation is a necessity because manually creating and dating many commits for a significant pattern would be extremely time-consuming. The tutorial will likely involve scripts or programs to generate commits and control their timestamps. You will need to utilize a scripting language like JavaScript with tools like Node.js to achieve the desired effect. This process is used to populate the graph and visually represent contributions over time.'
$d.Paragraphs(25).Range.Text = '*   Automation is crucial for generating numerous commits.'
$d.Paragraphs(26).Range.Text = '*   Scripts are typically used to create and date commits.'
$d.Paragraphs(27).Range.Text = '*   The focus is on creating the appearance of a long history.'
$d.Paragraphs(28).Range.Text = '*   Automated commits link with past dates.'
$d.Paragraphs(29).Range.Text = '*   This allows for the creation of complex contribution patterns.'
$d.Paragraphs(31).Range.Delete()

# --- Topic block 4 ---
$d.Paragraphs(31).Range.Text = 'Topic: Node.js Project Setup'
$d.Paragraphs(33).Range.Text = 'The tutorial uses Node.js and npm to set up a project, demonstrating that the technique uses JavaScript to manipulate the GitHub contribution graph. You''ll initialize a new project using ''npm init -y'' and this creates a package.json file. Then, you''ll create a JavaScript file, which will likely contain the logic for generating and dating the commits. This setup provides the necessary environment to run scripts that control the commit process. This approach is favored because it provides flexibility, and many readily available modules exist to make tasks like file manipulation and system calls easier, thus streamlining the commit process.'
$d.Paragraphs(35).Range.Text = '*   Node.js and npm are used to set up the project.'
$d.Paragraphs(36).Range.Text = '*   ''npm init -y'' initializes the project and creates a package.json file.'
$d.Paragraphs(37).Range.Text = '*   The JavaScript file holds the commit generation logic.'
$d.Paragraphs(38).Range.Text = '*   This environment facilitates the creation of the manipulation scripts.'
$d.Paragraphs(39).Range.Text = '*   It provides access to useful modules for file handling.'
$d.Paragraphs(41).Range.Delete()

# --- Topic block 5 ---
$d.Paragraphs(41).Range.Text = 'Topic: Using JSON Files for Commit Data'
$d.Paragraphs(43).Range.Text = 'The tutorial employs JSON files to store commit timestamp data. This file likely acts as a data source, containing information about the dates and times to which commits should be assigned. The JSON format is well-suited for storing structured data, making it easy to manage and parse the commit information. The JSON file enables you to define the pattern of commits you want to create and is accessible to be read by the JavaScript code. JSON is commonly used for configuration files and is a straightforward way to define the timestamps for the commits, making it easier to design different visual patterns on the GitHub contribution graph.'
$d.Paragraphs(45).Range.Text = '*   JSON files store the commit timestamp data.'
$d.Paragraphs(46).Range.Text = '*   The files define the desired commit pattern.'
$d.Paragraphs(47).Range.Text = '*   JSON is suitable for storing structured data.'
$d.Paragraphs(48).Range.Text = '*   The file is read by the JavaScript code.'
$d.Paragraphs(49).Range.Text = '*   It makes it easy to plan out different patterns.'
$d.Paragraphs(50).Range.Text = ''
$d.Paragraphs(51).Range.Delete()

# --- New Topic block 6 (insert) ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$insertAt = $d.Paragraphs.Count
$d.Paragraphs($insertAt).Range.Text = 'Topic: Using NPM Modules'
$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = ''
$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = 'The video involves the use of various npm modules to facilitate the GitHub contribution graph manipulation process. These modules provide specific functionalities, such as interacting with the file system, creating commits, and manipulating timestamps. The tutorial explains how these modules are installed using ''npm install'' and likely explains their individual roles within the JavaScript code. Understanding the use of these modules is crucial for replicating the technique. Each module serves a specific purpose in the code. The use of npm modules demonstrates an effective software development practice, which reuses functionality developed by others.'
$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = ''
$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = '*   npm modules provide specific functionalities.'
$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = '*   The modules are installed via ''npm install''.'
$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = '*   These modules support file management and timestamp handling.'
$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = '*   Knowing these modules is essential to use this technique.'
$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = '*   Modules represent efficient code reuse.'
$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = ''
